# "Generate Report for Archive"
# The localization status for this file moved on from "Ready for handoff"
# to "In Translation" on every sheet that tracks it (the Overview rollup
# columns for zh-cn / de-de, and each language sheet's own Status column).
# Shrinking that text also lets Excel's column autosize pull the Status
# columns in narrower, so we tighten those columns to match.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Narrow "custom" width Excel settles on once the longer status text is
# replaced by the shorter one (re-fit to content).
$newStatusColWidth = 12.58

# --- Overview sheet: zh-cn (E) and de-de (F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth
